# Updates cryptos list price/volume figures (and a couple of reordered coin
# rows) to match the latest scrape, per the "Updated cryptos list" commit.
#
# Note: some Price values look numeric (e.g. "1.00", "572.33"). Excel's
# COM layer auto-converts such strings to real numbers when assigned via
# .Value, which would silently drop formatting like trailing zeros
# ("1.00" -> 1). Setting NumberFormat to "@" (Text) first forces the
# assignment to keep the literal string, matching the original inline
# string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.302.40"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "3.426.78"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.33"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.49"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.426.91"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.555"
$ws.Range("E9").Value = "  -8.22%  "
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.121"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.425"
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("D13").Value = "4.010.39"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.27"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("E16").Value = "  -6.59%  "
$ws.Range("D17").Value = "64.333.65"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "3.455.14"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("E19").Value = "  -3.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.63"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.40"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.90"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.62"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.520"
$ws.Range("E25").Value = "  -5.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.08"
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.04"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.12"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  -5.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.06"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.862"
$ws.Range("E37").Value = "  +11.60%  "
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0734"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("D40").Value = "2.818.90"
$ws.Range("E40").Value = "  -3.86%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.57"
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.24"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.86"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.34"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.47"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0305"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "343.92"
$ws.Range("E47").Value = "  +8.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +5.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.07"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.33"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("E51").Value = "  -4.18%  "
